# Applies scheduled-runner updates to the Asura_Profits workbook.
# Updates pricing/profit columns (H-N) across the ALC, ARM, BSM, CRP, CUL,
# GSM, LTW and WVR sheets to reflect refreshed market data.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 51 (context G51=5486)
$ws.Range("H51").Value = 5000
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4516
# Row 107 (context G107=27766)
$ws.Range("H107").Value = 557.9524
$ws.Range("I107").Value = 321.5
$ws.Range("J107").Value = 772.9091
$ws.Range("K107").Value = 321.5
$ws.Range("L107").Value = 772.9091
$ws.Range("M107").Value = 1598.5
$ws.Range("N107").Value = -4612.9091
# Row 112 (context G112=27960)
$ws.Range("H112").Value = 2267.7407
$ws.Range("J112").Value = 2267.7407
$ws.Range("L112").Value = 6803.222099999999
$ws.Range("N112").Value = -9019.222099999999
# Row 116 (context G116=27778)
$ws.Range("H116").Value = 8002205.5
$ws.Range("I116").Value = 12501926
$ws.Range("J116").Value = 2701.6667
$ws.Range("K116").Value = 12501926
$ws.Range("L116").Value = 2701.6667
$ws.Range("M116").Value = -12498484
$ws.Range("N116").Value = -9585.6667
# Row 129 (context G129=36115)
$ws.Range("H129").Value = 1065
$ws.Range("J129").Value = 1288.1
$ws.Range("L129").Value = 3864.3
$ws.Range("N129").Value = -13864.3
# Row 132 (context G132=44049)
$ws.Range("H132").Value = 1820.2982
$ws.Range("I132").Value = 1380.5106
$ws.Range("J132").Value = 3887.3
$ws.Range("K132").Value = 4141.531800000001
$ws.Range("L132").Value = 11661.9
$ws.Range("M132").Value = -1611.531800000001
$ws.Range("N132").Value = -16721.9
# Row 135 (context G135=44047)
$ws.Range("H135").Value = 917.2353
$ws.Range("I135").Value = 854.72
$ws.Range("J135").Value = 1090.8889
$ws.Range("K135").Value = 7692.48
$ws.Range("L135").Value = 9818.0001
$ws.Range("M135").Value = -5157.48
$ws.Range("N135").Value = -14888.0001
# Row 137 (context G137=44013)
$ws.Range("H137").Value = 1336.1333
$ws.Range("I137").Value = 1215.875
$ws.Range("K137").Value = 3647.625
$ws.Range("M137").Value = -1097.625
# Row 138 (context G138=44169)
$ws.Range("H138").Value = 2364.679
$ws.Range("I138").Value = 1255.8223
$ws.Range("J138").Value = 3750.75
$ws.Range("K138").Value = 3767.4669
$ws.Range("L138").Value = 11252.25
$ws.Range("M138").Value = 1372.5331
$ws.Range("N138").Value = -21532.25
# Row 141 (context G141=44161)
$ws.Range("H141").Value = 5113.0625
$ws.Range("I141").Value = 2049.12
$ws.Range("J141").Value = 16055.714
$ws.Range("K141").Value = 6147.36
$ws.Range("L141").Value = 48167.142
$ws.Range("M141").Value = -967.3599999999997
$ws.Range("N141").Value = -58527.142

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (context G32=44147)
$ws.Range("H32").Value = 13008.5205
$ws.Range("I32").Value = 12771.613
$ws.Range("J32").Value = 15614.5
$ws.Range("K32").Value = 12771.613
$ws.Range("L32").Value = 15614.5
$ws.Range("M32").Value = -12484.613
$ws.Range("N32").Value = -16188.5
# Row 74 (context G74=44000)
$ws.Range("H74").Value = 1009.34283
$ws.Range("I74").Value = 929.0357
$ws.Range("K74").Value = 929.0357
$ws.Range("M74").Value = -55.03570000000002
# Row 77 (context G77=44000)
$ws.Range("H77").Value = 1009.34283
$ws.Range("I77").Value = 929.0357
$ws.Range("K77").Value = 4645.1785
$ws.Range("M77").Value = -277.1785
# Row 132 (context G132=43997)
$ws.Range("H132").Value = 4520.4565
$ws.Range("I132").Value = 5249.5713
$ws.Range("K132").Value = 15748.7139
$ws.Range("M132").Value = -13218.7139
# Row 135 (context G135=42016)
$ws.Range("H135").Value = 67237.664
$ws.Range("J135").Value = 67237.664
$ws.Range("L135").Value = 67237.664
$ws.Range("N135").Value = -77377.664

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134 (context G134=43998)
$ws.Range("H134").Value = 1669.3175
$ws.Range("I134").Value = 1424.1041
$ws.Range("J134").Value = 2454
$ws.Range("K134").Value = 4272.3123
$ws.Range("L134").Value = 7362
$ws.Range("M134").Value = -1737.3123
$ws.Range("N134").Value = -12432

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16 (context G16=27691)
$ws.Range("H16").Value = 1253.9333
$ws.Range("I16").Value = 1140.9
$ws.Range("J16").Value = 1480
$ws.Range("K16").Value = 1140.9
$ws.Range("L16").Value = 1480
$ws.Range("M16").Value = -853.9000000000001
$ws.Range("N16").Value = -2054
# Row 31 (context G31=44023)
$ws.Range("H31").Value = 2189.7878
$ws.Range("I31").Value = 1338.2
$ws.Range("J31").Value = 3499.923
$ws.Range("K31").Value = 1338.2
$ws.Range("L31").Value = 3499.923
$ws.Range("M31").Value = -1043.2
$ws.Range("N31").Value = -4089.923
# Row 34 (context G34=44023)
$ws.Range("H34").Value = 2189.7878
$ws.Range("I34").Value = 1338.2
$ws.Range("J34").Value = 3499.923
$ws.Range("K34").Value = 1338.2
$ws.Range("L34").Value = 3499.923
$ws.Range("M34").Value = -1136.2
$ws.Range("N34").Value = -3903.923
# Row 58 (context G58=44021)
$ws.Range("H58").Value = 1544943
$ws.Range("I58").Value = 1951073.4
$ws.Range("K58").Value = 1951073.4
$ws.Range("M58").Value = -1950870.4
# Row 113 (context G113=27691)
$ws.Range("H113").Value = 1253.9333
$ws.Range("I113").Value = 1140.9
$ws.Range("J113").Value = 1480
$ws.Range("K113").Value = 1140.9
$ws.Range("L113").Value = 1480
$ws.Range("M113").Value = 1029.1
$ws.Range("N113").Value = -5820
# Row 134 (context G134=44020)
$ws.Range("H134").Value = 1645.2894
$ws.Range("I134").Value = 1257.7273
$ws.Range("J134").Value = 2178.1875
$ws.Range("K134").Value = 3773.1819
$ws.Range("L134").Value = 6534.5625
$ws.Range("M134").Value = -1238.1819
$ws.Range("N134").Value = -11604.5625
# Row 136 (context G136=44021)
$ws.Range("H136").Value = 1544943
$ws.Range("I136").Value = 1951073.4
$ws.Range("K136").Value = 5853220.199999999
$ws.Range("M136").Value = -5850670.199999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68 (context G68=12895)
$ws.Range("H68").Value = 1256
$ws.Range("I68").Value = 1070
$ws.Range("K68").Value = 3210
$ws.Range("M68").Value = -2399
# Row 71 (context G71=12895)
$ws.Range("H71").Value = 1256
$ws.Range("I71").Value = 1070
$ws.Range("K71").Value = 9630
$ws.Range("M71").Value = -5574
# Row 138 (context G138=44105)
$ws.Range("H138").Value = 2180.2632
$ws.Range("I138").Value = 868.25
$ws.Range("J138").Value = 4429.4287
$ws.Range("K138").Value = 2604.75
$ws.Range("L138").Value = 13288.2861
$ws.Range("M138").Value = 2535.25
$ws.Range("N138").Value = -23568.2861

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 107 (context G107=27802)
$ws.Range("H107").Value = 883.4783
$ws.Range("I107").Value = 861.5
$ws.Range("J107").Value = 933.7143
$ws.Range("K107").Value = 861.5
$ws.Range("L107").Value = 933.7143
$ws.Range("M107").Value = 1058.5
$ws.Range("N107").Value = -4773.7143

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61 (context G61=27740)
$ws.Range("H61").Value = 14347.529
$ws.Range("I61").Value = 18054.46
$ws.Range("J61").Value = 2300
$ws.Range("K61").Value = 18054.46
$ws.Range("L61").Value = 2300
$ws.Range("M61").Value = -17852.46
$ws.Range("N61").Value = -2704
# Row 113 (context G113=27740)
$ws.Range("H113").Value = 14347.529
$ws.Range("I113").Value = 18054.46
$ws.Range("J113").Value = 2300
$ws.Range("K113").Value = 18054.46
$ws.Range("L113").Value = 2300
$ws.Range("M113").Value = -15884.46
$ws.Range("N113").Value = -6640
# Row 136 (context G136=44060)
$ws.Range("H136").Value = 15076084
$ws.Range("I136").Value = 20409266
$ws.Range("J136").Value = 557975.7
$ws.Range("K136").Value = 61227798
$ws.Range("L136").Value = 1673927.1
$ws.Range("M136").Value = -61225248
$ws.Range("N136").Value = -1679027.1

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 126 (context G126=36210)
$ws.Range("H126").Value = 10999
$ws.Range("I126").Value = 10999
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 32997
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -30527
$ws.Range("N126").ClearContents()
# Row 132 (context G132=44029)
$ws.Range("H132").Value = 1880.9535
$ws.Range("I132").Value = 1178.3334
$ws.Range("J132").Value = 4199.6
$ws.Range("K132").Value = 3535.0002
$ws.Range("L132").Value = 12598.8
$ws.Range("M132").Value = -1005.0002
$ws.Range("N132").Value = -17658.8
# Row 136 (context G136=44031)
$ws.Range("H136").Value = 1709.6296
$ws.Range("I136").Value = 1550.6522
$ws.Range("K136").Value = 4651.9566
$ws.Range("M136").Value = -2101.9566
# Row 137 (context G137=42184)
$ws.Range("H137").Value = 38460.715
$ws.Range("J137").Value = 38460.715
$ws.Range("L137").Value = 38460.715
$ws.Range("N137").Value = -48660.715
